$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the old column D, now shifted) into new D:E
$ws.Columns("F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new quarter columns (D = period ending 2018-09-30, E = period ending 2018-06-30)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 3299000
$ws.Range("E8").Value = 3611000
$ws.Range("D9").Value = 343000
$ws.Range("E9").Value = 337000
$ws.Range("D10").Value = 2956000
$ws.Range("E10").Value = 3274000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 12000
$ws.Range("E14").Value = 3000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 3387000
$ws.Range("E17").Value = 3111000
$ws.Range("D18").Value = -88000
$ws.Range("E18").Value = 500000
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 144000
$ws.Range("E22").Value = 146000
$ws.Range("D23").Value = -232000
$ws.Range("E23").Value = 354000
$ws.Range("D24").Value = -21000
$ws.Range("E24").Value = 65000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -211000
$ws.Range("E26").Value = 289000
$ws.Range("D27").Value = -165000
$ws.Range("E27").Value = 278000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = -165000
$ws.Range("E33").Value = 278000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -165000
$ws.Range("E35").Value = 278000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 405000
$ws.Range("E41").Value = 571000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 7960000
$ws.Range("E43").Value = 7836000
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 48186000
$ws.Range("E47").Value = 48695000
$ws.Range("D48").Value = 15511000
$ws.Range("E48").Value = 15467000
$ws.Range("D49").Value = 665000
$ws.Range("E49").Value = 657000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 78316000
$ws.Range("E54").Value = 78725000
$ws.Range("D57").Value = 0
$ws.Range("E57").Value = 0
$ws.Range("D58").Value = 17000
$ws.Range("E58").Value = 149000
$ws.Range("D59").Value = 4225000
$ws.Range("E59").Value = 4506000
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 11359000
$ws.Range("E61").Value = 11311000
$ws.Range("D62").Value = 841000
$ws.Range("E62").Value = 911000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 59798000
$ws.Range("E66").Value = 59773000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 15773000
$ws.Range("E72").Value = 16790000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 18518000
$ws.Range("E76").Value = 18952000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -165000
$ws.Range("E81").Value = 278000
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 827000
$ws.Range("E89").Value = 1244000
$ws.Range("D91").Value = -264000
$ws.Range("E91").Value = -251000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -727000
$ws.Range("E94").Value = 520000
$ws.Range("D96").Value = -20000
$ws.Range("E96").Value = -20000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -260000
$ws.Range("E100").Value = -1647000
$ws.Range("D101").Value = -6000
$ws.Range("E101").Value = 1000
$ws.Range("D102").Value = -166000
$ws.Range("E102").Value = 118000
